$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header cell (G1) onto the new H1 header
# so the "Save" header picks up the same bold/border/alignment style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 8).Value = "Save"

# Values for the new "Save" column (rows 2-11)
$saveValues = @(0, 1, 0, 0, 0, 0, 0, 1, 1, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = 0
